# Scheduled-runner data refresh: overwrites literal market-board-derived
# price/profit columns (H..N) for specific Leve rows across the ALC, ARM,
# BSM, CRP, CUL, GSM, LTW and WVR sheets. Source data has no formulas -
# every touched cell is a plain literal value, so we just re-assign .Value
# (or ClearContents() where a cell is removed entirely, i.e. goes blank).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 3914.9119
$ws.Range("J17").Value = 4087.7188
$ws.Range("L17").Value = 12263.1564
$ws.Range("N17").Value = -12599.1564
# row 19
$ws.Range("H19").Value = 2847.0476
$ws.Range("I19").Value = 953.9
$ws.Range("K19").Value = 953.9
$ws.Range("M19").Value = -778.9
# row 31
$ws.Range("H31").Value = 1150.2
$ws.Range("I31").Value = 1150.2
$ws.Range("K31").Value = 3450.6
$ws.Range("M31").Value = -3220.6
# row 33
$ws.Range("H33").Value = 715.64703
$ws.Range("I33").Value = 749.75
$ws.Range("K33").Value = 749.75
$ws.Range("M33").Value = -520.75
# row 53
$ws.Range("H53").Value = 1026.875
$ws.Range("I53").Value = 381.85715
$ws.Range("K53").Value = 381.85715
$ws.Range("M53").Value = 255.14285
# row 64
$ws.Range("H64").Value = 7685.154
$ws.Range("I64").Value = 4999.6665
$ws.Range("K64").Value = 4999.6665
$ws.Range("M64").Value = -4751.6665
# row 67
$ws.Range("H67").Value = 7685.154
$ws.Range("I67").Value = 4999.6665
$ws.Range("K67").Value = 4999.6665
$ws.Range("M67").Value = -4141.6665
# row 74
$ws.Range("H74").Value = 7985
$ws.Range("J74").Value = 7985
$ws.Range("L74").Value = 7985
$ws.Range("N74").Value = -9857
# row 77
$ws.Range("H77").Value = 7985
$ws.Range("J77").Value = 7985
$ws.Range("L77").Value = 39925
$ws.Range("N77").Value = -49285
# row 112
$ws.Range("H112").Value = 11270.9
$ws.Range("J112").Value = 13876.125
$ws.Range("L112").Value = 41628.375
$ws.Range("N112").Value = -43844.375
# row 113
$ws.Range("H113").Value = 4419.737
$ws.Range("I113").Value = 1930.875
$ws.Range("K113").Value = 1930.875
$ws.Range("M113").Value = 1323.125
# row 138
$ws.Range("H138").Value = 3127.9841
$ws.Range("I138").Value = 1756
$ws.Range("J138").Value = 3450.804
$ws.Range("K138").Value = 5268
$ws.Range("L138").Value = 10352.412
$ws.Range("M138").Value = -128
$ws.Range("N138").Value = -20632.412

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 5330659
$ws.Range("I45").Value = 7194145.5
$ws.Range("J45").Value = 6413
$ws.Range("K45").Value = 7194145.5
$ws.Range("L45").Value = 6413
$ws.Range("M45").Value = -7193768.5
$ws.Range("N45").Value = -7167
# row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# row 122
$ws.Range("H122").Value = 497987.97
$ws.Range("I122").Value = 1602.7567
$ws.Range("K122").Value = 4808.2701
$ws.Range("M122").Value = -2358.2701
# row 132
$ws.Range("H132").Value = 3012
$ws.Range("I132").Value = 2147.6365
$ws.Range("K132").Value = 6442.9095
$ws.Range("M132").Value = -3912.9095

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 2895.2554
$ws.Range("I134").Value = 1259.2812
$ws.Range("J134").Value = 6385.3335
$ws.Range("K134").Value = 3777.8436
$ws.Range("L134").Value = 19156.0005
$ws.Range("M134").Value = -1242.8436
$ws.Range("N134").Value = -24226.0005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 32267.232
$ws.Range("J31").Value = 93922.2
$ws.Range("L31").Value = 93922.2
$ws.Range("N31").Value = -94512.2
# row 34
$ws.Range("H34").Value = 32267.232
$ws.Range("J34").Value = 93922.2
$ws.Range("L34").Value = 93922.2
$ws.Range("N34").Value = -94326.2
# row 86
$ws.Range("H86").Value = 8657.682000000001
$ws.Range("I86").Value = 6964.1816
$ws.Range("J86").Value = 10351.182
$ws.Range("K86").Value = 6964.1816
$ws.Range("L86").Value = 10351.182
$ws.Range("M86").Value = -5841.1816
$ws.Range("N86").Value = -12597.182
# row 89
$ws.Range("H89").Value = 8657.682000000001
$ws.Range("I89").Value = 6964.1816
$ws.Range("J89").Value = 10351.182
$ws.Range("K89").Value = 34820.908
$ws.Range("L89").Value = 51755.91
$ws.Range("M89").Value = -29204.908
$ws.Range("N89").Value = -62987.91

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 12
$ws.Range("H12").Value = 37132.293
$ws.Range("J12").Value = 112.833336
$ws.Range("L12").Value = 338.500008
$ws.Range("N12").Value = -684.500008
# row 23
$ws.Range("H23").Value = 146.23077
$ws.Range("I23").Value = 29.8
$ws.Range("J23").Value = 219
$ws.Range("K23").Value = 89.40000000000001
$ws.Range("L23").Value = 657
$ws.Range("M23").Value = 145.6
$ws.Range("N23").Value = -1127
# row 99
$ws.Range("H99").Value = 4995
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4995
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 14985
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -19477
# row 130
$ws.Range("H130").Value = 2621.5
$ws.Range("I130").Value = 865
$ws.Range("K130").Value = 2595
$ws.Range("M130").Value = 2425

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 43
$ws.Range("H43").Value = 13408.833
$ws.Range("I43").Value = 1472.3334
$ws.Range("K43").Value = 1472.3334
$ws.Range("M43").Value = -1321.3334
# row 126
$ws.Range("H126").Value = 5154704
$ws.Range("I126").Value = 3032963
$ws.Range("K126").Value = 9098889
$ws.Range("M126").Value = -9096419
# row 132
$ws.Range("H132").Value = 3603.138
$ws.Range("I132").Value = 3207.4666
$ws.Range("K132").Value = 9622.399800000001
$ws.Range("M132").Value = -7092.399800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3631.9473
$ws.Range("I7").Value = 1820.4
$ws.Range("J7").Value = 5644.778
$ws.Range("K7").Value = 1820.4
$ws.Range("L7").Value = 5644.778
$ws.Range("M7").Value = -1708.4
$ws.Range("N7").Value = -5868.778
# row 46
$ws.Range("H46").Value = 4099.6665
$ws.Range("I46").Value = 2483
$ws.Range("J46").Value = 7333
$ws.Range("K46").Value = 2483
$ws.Range("L46").Value = 7333
$ws.Range("M46").Value = -2295
$ws.Range("N46").Value = -7709
# row 68
$ws.Range("H68").Value = 2725
$ws.Range("J68").Value = 3255.8
$ws.Range("L68").Value = 3255.8
$ws.Range("N68").Value = -4753.8
# row 71
$ws.Range("H71").Value = 2725
$ws.Range("J71").Value = 3255.8
$ws.Range("L71").Value = 16279
$ws.Range("N71").Value = -23767
# row 126
$ws.Range("H126").Value = 3631.9473
$ws.Range("I126").Value = 1820.4
$ws.Range("J126").Value = 5644.778
$ws.Range("K126").Value = 5461.200000000001
$ws.Range("L126").Value = 16934.334
$ws.Range("M126").Value = -2991.200000000001
$ws.Range("N126").Value = -21874.334
# row 132
$ws.Range("H132").Value = 5136.755
$ws.Range("I132").Value = 4447
$ws.Range("J132").Value = 6435.1177
$ws.Range("K132").Value = 13341
$ws.Range("L132").Value = 19305.3531
$ws.Range("M132").Value = -10811
$ws.Range("N132").Value = -24365.3531

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 107
$ws.Range("H107").Value = 43480676
$ws.Range("I107").Value = 50002684
$ws.Range("J107").Value = 637.6667
$ws.Range("K107").Value = 150008052
$ws.Range("L107").Value = 1913.0001
$ws.Range("M107").Value = -150006132
$ws.Range("N107").Value = -5753.0001
